$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

# Update the Week 15 distribution simulation strings (B2:E2)
$ws.Range("B2").Value = "NIG(1.0436408138409994, 0.8124501614504589, 0.8808321659619507, 3.1274472362764705)"
$ws.Range("C2").Value = "JSU(-0.8222644097515871, 1.05202803488568, 4.744755414132969, 4.363954381562556)"
$ws.Range("D2").Value = "NIG(0.544995928668363, 0.31243505445076236, 2.6436849570409198, 2.788357328496261)"
$ws.Range("E2").Value = "JSU(-1.0711820559001244, 1.1642134631721452, 3.297954423379929, 5.656121892825822)"
